# Added internal gains in heating shifting, revised units of measure in launcher_shift2
# Update computed results on row 13 ("(M)+ECS" combination) of sheet "4F".

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("4F")

$row = 13

$ws.Cells.Item($row, 21).Value  = 10                     # U13  - PV [kWp]
$ws.Cells.Item($row, 23).Value  = 15000                  # W13  - Investissement PV [EUR]

$ws.Cells.Item($row, 33).Value  = 7.727832487873169      # AG13 - Capacite d'acces kW
$ws.Cells.Item($row, 34).Value  = 11559.42574525542       # AH13 - Total [kWh]
$ws.Cells.Item($row, 35).Value  = 10566.36352937754       # AI13 - Electricite Produite [kWh]
$ws.Cells.Item($row, 36).Value  = 3521.330594287391        # AJ13 - Electricite autoconsommee [kWh]
$ws.Cells.Item($row, 37).Value  = 7045.032935090147        # AK13 - Electricite injectee [kWh]
$ws.Cells.Item($row, 38).Value  = 8038.095150968032        # AL13 - Electricite achetee [kWh]
$ws.Cells.Item($row, 39).Value  = 607.1753285511822         # AM13 - Heure Talon [kWh]
$ws.Cells.Item($row, 40).Value  = 1482.44036143576          # AN13 - Heure creuse [kWh]
$ws.Cells.Item($row, 41).Value  = 4116.063933945034         # AO13 - Heure pleine [kWh]
$ws.Cells.Item($row, 42).Value  = 1832.415527036056         # AP13 - Heure pointe [kWh]
$ws.Cells.Item($row, 43).Value  = -0.6345958244451267        # AQ13 - Variation Heure Talon
$ws.Cells.Item($row, 44).Value  = -8.42685369907935           # AR13 - Variation Heure creuse
$ws.Cells.Item($row, 45).Value  = -21.29685747924691          # AS13 - Variation Heure pleine
$ws.Cells.Item($row, 46).Value  = -0.1045435911062923         # AT13 - Variation Heure pointe
$ws.Cells.Item($row, 47).Value  = 0.3046285059387767           # AU13 - Taux autosuff [%]
$ws.Cells.Item($row, 48).Value  = 0.3332585126848113           # AV13 - Taux autocons [%]

$ws.Cells.Item($row, 50).Value  = 281.8013174036062            # AX13
$ws.Cells.Item($row, 51).Value  = 2286.96504063424              # AY13

$ws.Cells.Item($row, 53).Value  = -2005.163723230634            # BA13
$ws.Cells.Item($row, 54).Value  = 0.2970564980271703             # BB13

$ws.Cells.Item($row, 56).Value  = -2185.954915917501             # BD13
$ws.Cells.Item($row, 57).Value  = -0.1457303277278334            # BE13
